$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, [string]$value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ---- Row 2: Oct 3 2020 vs Delhi Capitals (Sharjah) ----
$ws.Range("A2").Value = " Oct 3 2020"
$ws.Range("B2").Value = " Sharjah"
$ws.Range("C2").Value = "Capitals won by 18 runs"
$ws.Range("D2").Value = "Kolkata Knight Riders"
$ws.Range("E2").Value = "Delhi Capitals"
$ws.Range("F2").Value = "Shivam Mavi "
Set-TextCell $ws.Range("G2") "1"
Set-TextCell $ws.Range("H2") "3"
Set-TextCell $ws.Range("I2") "0"
Set-TextCell $ws.Range("J2") "0"
Set-TextCell $ws.Range("K2") "33.33"

# ---- Row 3: Oct 7 2020 vs Chennai Super Kings (Abu Dhabi) ----
$ws.Range("A3").Value = " Oct 7 2020"
$ws.Range("B3").Value = " Abu Dhabi"
$ws.Range("C3").Value = "KKR won by 10 runs"
$ws.Range("D3").Value = "Kolkata Knight Riders"
$ws.Range("E3").Value = "Chennai Super Kings"
$ws.Range("F3").Value = "Shivam Mavi "
Set-TextCell $ws.Range("G3") "0"
Set-TextCell $ws.Range("H3") "1"
Set-TextCell $ws.Range("I3") "0"
Set-TextCell $ws.Range("J3") "0"
Set-TextCell $ws.Range("K3") "0.00"

# ---- Row 4: Sep 23 2020 vs Mumbai Indians (Abu Dhabi) ----
$ws.Range("A4").Value = " Sep 23 2020"
$ws.Range("B4").Value = " Abu Dhabi"
$ws.Range("C4").Value = "Mumbai won by 49 runs"
$ws.Range("D4").Value = "Kolkata Knight Riders"
$ws.Range("E4").Value = "Mumbai Indians"
$ws.Range("F4").Value = "Shivam Mavi "
Set-TextCell $ws.Range("G4") "9"
Set-TextCell $ws.Range("H4") "10"
Set-TextCell $ws.Range("I4") "1"
Set-TextCell $ws.Range("J4") "0"
Set-TextCell $ws.Range("K4") "90.00"
